$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain plain text so numeric-looking
# strings (e.g. "1.000", "0.9991") are not silently coerced into numbers
# and lose precision/trailing zeros.
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @{
    'D2' = '30.116.16'
    'E2' = '  +10.12%  '
    'D3' = '1.871.37'
    'E3' = '  +7.00%  '
    'D4' = '1.000'
    'E4' = '  +0.12%  '
    'D5' = '250.06'
    'E5' = '  +3.38%  '
    'D6' = '1.000'
    'E6' = '  +0.09%  '
    'D7' = '0.4955'
    'E7' = '  +2.94%  '
    'E8' = '  +8.86%  '
    'D9' = '0.2839'
    'E9' = '  +8.37%  '
    'D10' = '0.06527'
    'E10' = '  +5.79%  '
    'D11' = '1.867.78'
    'E11' = '  +6.77%  '
    'E12' = '  +5.56%  '
    'D13' = '0.07192'
    'E13' = '  +3.58%  '
    'D14' = '0.6608'
    'E14' = '  +9.38%  '
    'D15' = '85.28'
    'E15' = '  +10.43%  '
    'D16' = '4.800'
    'E16' = '  +7.35%  '
    'D17' = '30.077.51'
    'E17' = '  +10.05%  '
    'D18' = '0.9991'
    'E18' = '  +0.05%  '
    'D19' = '0.000007488'
    'E19' = '  +5.66%  '
    'D20' = '12.67'
    'E20' = '  +10.46%  '
    'D21' = '1.000'
    'E21' = '  +0.08%  '
    'D22' = '2.110.37'
    'E22' = '  +7.59%  '
    'D23' = '4.716'
    'E23' = '  +6.03%  '
    'D24' = '5.508'
    'E24' = '  +7.71%  '
    'D25' = '8.997'
    'E25' = '  +6.58%  '
    'D26' = '144.39'
    'E26' = '  +1.62%  '
    'D27' = '134.54'
    'E27' = '  +25.04%  '
    'D28' = '16.75'
    'E28' = '  +9.67%  '
    'D29' = '1.941'
    'E29' = '  +5.56%  '
    'D30' = '1.397'
    'E30' = '  -0.89%  '
    'E31' = '  +7.42%  '
    'D32' = '0.08597'
    'E32' = '  +7.82%  '
    'D33' = '3.884'
    'E33' = '  +5.73%  '
    'D34' = '0.05054'
    'E34' = '  +8.03%  '
    'D35' = '1.130'
    'E35' = '  +11.32%  '
    'D36' = '0.6830'
    'E36' = '  +10.44%  '
    'D37' = '2.685'
    'E37' = '  +3.38%  '
    'D38' = '2.330'
    'E38' = '  +15.61%  '
    'D39' = '2.732'
    'E39' = '  +6.95%  '
    'D40' = '0.9562'
    'E40' = '  +3.52%  '
    'D41' = '0.01632'
    'E41' = '  +9.28%  '
    'D42' = '6.142'
    'E42' = '  +7.67%  '
    'D43' = '1.001'
    'E43' = '  +0.21%  '
    'D44' = '102.68'
    'E44' = '  +3.02%  '
    'D45' = '0.4167'
    'E45' = '  +8.38%  '
    'D46' = '7.390'
    'E46' = '  +7.38%  '
    'E47' = '  +8.16%  '
    'D48' = '0.05631'
    'E48' = '  +4.98%  '
    'B49' = 'EnergySwap'
    'C49' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D49' = '8.292'
    'E49' = '  +6.13%  '
    'B50' = 'Elrond'
    'C50' = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
    'D50' = '32.34'
    'E50' = '  +8.38%  '
    'B51' = 'Decentraland'
    'C51' = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    'D51' = '0.3710'
    'E51' = '  +9.77%  '
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# Restore the default (General) style on the Price/Volume columns so the
# saved workbook does not carry a stray text-format style on these cells.
$ws.Range("D2:E51").Style = "Normal"

